$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.409.55'
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.850.47'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.11%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.66%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6291'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.84%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07672'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.32%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2976'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.11%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.46'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.40%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.960.43'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.62%  '

# Row 12
$ws.Range("E12").Value = '  +1.00%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.006'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.83%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6889'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.28%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.18'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.42%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009993'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.17%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.177.76'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.14%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.189'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.93%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.552.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.63%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '232.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.98%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.47%  '

# Row 22
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9997'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.05%  '

# Row 23
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.661'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.85%  '

# Row 24
$ws.Range("E24").Value = '  -0.08%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.82'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.16%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1393'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.20%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.486'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.56%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.06%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.478'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.06%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05792'
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = '  -1.15%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.128'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.35%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.022'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.23%  '

# Row 34
$ws.Range("E34").Value = '  +0.74%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.162'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.89%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7242'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.45%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.587'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.47%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.248.71'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.83%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.797'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.15%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01807'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.05%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9095'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.32%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.094'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.10%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.106.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.31%  '

# Row 44
$ws.Range("E44").Value = '  -0.05%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '67.91'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.95%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.72'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.59%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.283'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.84%  '

# Row 48
$ws.Range("E48").Value = '  -0.09%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.208'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.53%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4032'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.53%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.699'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.43%  '
